# Microsite Education Script completed
# Appends the new interview-history rows captured for sprint 165 (AMSIN/BETA/AMS)
# and sprint 166 (AMSIN/BETA/AMS) to the respective worksheets.

$wb = $excel.ActiveWorkbook

function Add-HistoryRow($ws, $row, $runDate, $runTime, $sprintName, $total, $pass, $fail, $timeTaken) {
    # Column A: Run Date - stored as literal text (matches existing rows),
    # force text format first so the date-like string isn't auto-converted
    # into a serial date number.
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $runDate

    # Column B: Run Time - numeric serial date/time value.
    $cellB = $ws.Cells.Item($row, 2)
    $cellB.NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $cellB.Value = $runTime

    # Column C: Sprint Name - literal text.
    $cellC = $ws.Cells.Item($row, 3)
    $cellC.NumberFormat = "@"
    $cellC.Value = $sprintName

    # Columns D-G: numeric counts / duration.
    $ws.Cells.Item($row, 4).Value = $total
    $ws.Cells.Item($row, 5).Value = $pass
    $ws.Cells.Item($row, 6).Value = $fail
    $ws.Cells.Item($row, 7).Value = $timeTaken
}

# --- AMSIN sheet: append rows 42-46 (sprint 165 + 166 runs) ---
$wsAmsin = $wb.Worksheets.Item("AMSIN")
Add-HistoryRow $wsAmsin 42 "2022-08-02" 44775.66775893518 "165_fstcycle"  105 105 0 2.93
Add-HistoryRow $wsAmsin 43 "2022-08-03" 44776.69424932871 "165_scndcycle" 105 105 0 3
Add-HistoryRow $wsAmsin 44 "2022-08-04" 44777.39674225695 "165_finalrun"  105 105 0 2.83
Add-HistoryRow $wsAmsin 45 "2022-08-22" 44795.67676381944 "166fstcycle"   105 100 5 3.8
Add-HistoryRow $wsAmsin 46 "2022-08-23" 44796.9109019213  "166cyclescnd"  105 105 0 2.84

# --- BETA sheet: append rows 22-23 (sprint 165 + 166 beta runs) ---
$wsBeta = $wb.Worksheets.Item("BETA")
Add-HistoryRow $wsBeta 22 "2022-08-04" 44777.56998666667 "165beta"  105 105 0 2.68
Add-HistoryRow $wsBeta 23 "2022-08-24" 44797.55008375    "166_beta" 105 105 0 2.71

# --- AMS sheet: append rows 21-22 (sprint 165 + 166 live runs) ---
$wsAms = $wb.Worksheets.Item("AMS")
Add-HistoryRow $wsAms 21 "2022-08-04" 44777.82335487269 "165_live" 105 105 0 2.74
Add-HistoryRow $wsAms 22 "2022-08-24" 44797.92717719739 "166_live" 105 105 0 2.72
